$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H69").Value = 3980
$ws.Range("J69").Value = 3980
$ws.Range("L69").Value = 11940
$ws.Range("N69").Value = -13688
$ws.Range("H72").Value = 3980
$ws.Range("J72").Value = 3980
$ws.Range("L72").Value = 35820
$ws.Range("N72").Value = -44556
$ws.Range("H107").Value = 2695.0386
$ws.Range("I107").Value = 1762.409
$ws.Range("K107").Value = 1762.409
$ws.Range("M107").Value = 157.5909999999999
$ws.Range("H112").Value = 1855.3125
$ws.Range("I112").Value = 452.4
$ws.Range("J112").Value = 2493
$ws.Range("K112").Value = 1357.2
$ws.Range("L112").Value = 7479
$ws.Range("M112").Value = -249.1999999999998
$ws.Range("N112").Value = -9695
$ws.Range("H129").Value = 851.5965
$ws.Range("I129").Value = 694
$ws.Range("K129").Value = 2082
$ws.Range("M129").Value = 2918
$ws.Range("H135").Value = 29412606
$ws.Range("I135").Value = 556.12
$ws.Range("J135").Value = 111112744
$ws.Range("K135").Value = 5005.08
$ws.Range("L135").Value = 1000014696
$ws.Range("M135").Value = -2470.08
$ws.Range("N135").Value = -1000019766
$ws.Range("H138").Value = 2346.97
$ws.Range("I138").Value = 1512.7693
$ws.Range("J138").Value = 2471.6206
$ws.Range("K138").Value = 4538.3079
$ws.Range("L138").Value = 7414.861800000001
$ws.Range("M138").Value = 601.6921000000002
$ws.Range("N138").Value = -17694.8618
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H22").Value = 782.75
$ws.Range("I22").Value = 782.75
$ws.Range("K22").Value = 782.75
$ws.Range("M22").Value = -483.75
$ws.Range("H32").Value = 9779.812
$ws.Range("I32").Value = 7225.9595
$ws.Range("K32").Value = 7225.9595
$ws.Range("M32").Value = -6938.9595
$ws.Range("H61").Value = 76924580
$ws.Range("I61").Value = 111112220
$ws.Range("J61").Value = 2374.5
$ws.Range("K61").Value = 111112220
$ws.Range("L61").Value = 2374.5
$ws.Range("M61").Value = -111112008
$ws.Range("N61").Value = -2798.5
$ws.Range("H63").Value = 2433.8438
$ws.Range("I63").Value = 2396.6
$ws.Range("J63").Value = 2992.5
$ws.Range("K63").Value = 2396.6
$ws.Range("L63").Value = 2992.5
$ws.Range("M63").Value = -1710.6
$ws.Range("N63").Value = -4364.5
$ws.Range("H66").Value = 2433.8438
$ws.Range("I66").Value = 2396.6
$ws.Range("J66").Value = 2992.5
$ws.Range("K66").Value = 11983
$ws.Range("L66").Value = 14962.5
$ws.Range("M66").Value = -8551
$ws.Range("N66").Value = -21826.5
$ws.Range("H74").Value = 4246.6665
$ws.Range("I74").Value = 3726.6667
$ws.Range("J74").Value = 4766.6665
$ws.Range("K74").Value = 3726.6667
$ws.Range("L74").Value = 4766.6665
$ws.Range("M74").Value = -2852.6667
$ws.Range("N74").Value = -6514.6665
$ws.Range("H77").Value = 4246.6665
$ws.Range("I77").Value = 3726.6667
$ws.Range("J77").Value = 4766.6665
$ws.Range("K77").Value = 18633.3335
$ws.Range("L77").Value = 23833.3325
$ws.Range("M77").Value = -14265.3335
$ws.Range("N77").Value = -32569.3325
$ws.Range("H97").Value = 565.5
$ws.Range("I97").Value = 542.3182
$ws.Range("K97").Value = 542.3182
$ws.Range("M97").Value = -46.31820000000005
$ws.Range("H103").Value = 43840.5
$ws.Range("J103").Value = 43840.5
$ws.Range("L103").Value = 43840.5
$ws.Range("N103").Value = -46184.5
$ws.Range("H136").Value = 76924580
$ws.Range("I136").Value = 111112220
$ws.Range("J136").Value = 2374.5
$ws.Range("K136").Value = 333336660
$ws.Range("L136").Value = 7123.5
$ws.Range("M136").Value = -333334110
$ws.Range("N136").Value = -12223.5
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H82").Value = 15655.846
$ws.Range("I82").Value = 2608.8572
$ws.Range("J82").Value = 30877.334
$ws.Range("K82").Value = 2608.8572
$ws.Range("L82").Value = 30877.334
$ws.Range("M82").Value = -2225.8572
$ws.Range("N82").Value = -31643.334
$ws.Range("H85").Value = 15655.846
$ws.Range("I85").Value = 2608.8572
$ws.Range("J85").Value = 30877.334
$ws.Range("K85").Value = 2608.8572
$ws.Range("L85").Value = 30877.334
$ws.Range("M85").Value = -1282.8572
$ws.Range("N85").Value = -33529.334
$ws.Range("H94").Value = 41667970
$ws.Range("J94").Value = 1900
$ws.Range("L94").Value = 1900
$ws.Range("N94").Value = -2802
$ws.Range("H105").Value = 111123350
$ws.Range("I105").Value = 111123350
$ws.Range("J105").Value = 0
$ws.Range("K105").Value = 111123350
$ws.Range("L105").Value = 0
$ws.Range("M105").ClearContents()
$ws.Range("N105").Value = -111121603
$ws.Range("H107").Value = 992.6818
$ws.Range("I107").Value = 863.8421
$ws.Range("K107").Value = 863.8421
$ws.Range("M107").Value = 1056.1579
$ws.Range("H134").Value = 18028.5
$ws.Range("I134").Value = 1543
$ws.Range("K134").Value = 4629
$ws.Range("M134").Value = -2094
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 5590.6
$ws.Range("I58").Value = 1153.25
$ws.Range("J58").Value = 8548.833000000001
$ws.Range("K58").Value = 1153.25
$ws.Range("L58").Value = 8548.833000000001
$ws.Range("M58").Value = -950.25
$ws.Range("N58").Value = -8954.833000000001
$ws.Range("H74").Value = 27916.666
$ws.Range("I74").Value = 17250
$ws.Range("J74").Value = 33250
$ws.Range("K74").Value = 17250
$ws.Range("L74").Value = 33250
$ws.Range("M74").Value = -16376
$ws.Range("N74").Value = -34998
$ws.Range("H77").Value = 27916.666
$ws.Range("I77").Value = 17250
$ws.Range("J77").Value = 33250
$ws.Range("K77").Value = 51750
$ws.Range("L77").Value = 99750
$ws.Range("M77").Value = -47382
$ws.Range("N77").Value = -108486
$ws.Range("H105").Value = 796.44446
$ws.Range("I105").Value = 771
$ws.Range("K105").Value = 771
$ws.Range("M105").Value = 976
$ws.Range("H132").Value = 3016.2307
$ws.Range("I132").Value = 2856.889
$ws.Range("J132").Value = 3374.75
$ws.Range("K132").Value = 8570.667000000001
$ws.Range("L132").Value = 10124.25
$ws.Range("M132").Value = -6040.667000000001
$ws.Range("N132").Value = -15184.25
$ws.Range("H134").Value = 21740866
$ws.Range("I134").Value = 1744.8125
$ws.Range("K134").Value = 5234.4375
$ws.Range("M134").Value = -2699.4375
$ws.Range("H136").Value = 5590.6
$ws.Range("I136").Value = 1153.25
$ws.Range("J136").Value = 8548.833000000001
$ws.Range("K136").Value = 3459.75
$ws.Range("L136").Value = 25646.499
$ws.Range("M136").Value = -909.75
$ws.Range("N136").Value = -30746.499
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H88").Value = 6047.5713
$ws.Range("J88").Value = 6047.5713
$ws.Range("L88").Value = 18142.7139
$ws.Range("N88").Value = -18998.7139
$ws.Range("H91").Value = 6047.5713
$ws.Range("J91").Value = 6047.5713
$ws.Range("L91").Value = 18142.7139
$ws.Range("N91").Value = -21106.7139
$ws.Range("H117").Value = 766.7
$ws.Range("I117").Value = 626.1111
$ws.Range("J117").Value = 2032
$ws.Range("K117").Value = 1878.3333
$ws.Range("L117").Value = 6096
$ws.Range("M117").Value = 1563.6667
$ws.Range("N117").Value = -12980
$ws.Range("H131").Value = 17268128
$ws.Range("J131").Value = 35079.797
$ws.Range("L131").Value = 105239.391
$ws.Range("N131").Value = -115319.391
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 1011.25
$ws.Range("I122").Value = 1011.25
$ws.Range("K122").Value = 3033.75
$ws.Range("M122").Value = -583.75
$ws.Range("H132").Value = 4139.143
$ws.Range("I132").Value = 3957.9285
$ws.Range("J132").Value = 4501.5713
$ws.Range("K132").Value = 11873.7855
$ws.Range("L132").Value = 13504.7139
$ws.Range("M132").Value = -9343.7855
$ws.Range("N132").Value = -18564.7139
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 4306.923
$ws.Range("I46").Value = 997.75
$ws.Range("J46").Value = 5777.6665
$ws.Range("K46").Value = 997.75
$ws.Range("L46").Value = 5777.6665
$ws.Range("M46").Value = -809.75
$ws.Range("N46").Value = -6153.6665
$ws.Range("H132").Value = 3022.1738
$ws.Range("I132").Value = 4266.3335
$ws.Range("K132").Value = 12799.0005
$ws.Range("M132").Value = -10269.0005
$ws.Range("H136").Value = 2360
$ws.Range("I136").Value = 1933.3334
$ws.Range("K136").Value = 5800.0002
$ws.Range("M136").Value = -3250.0002
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 3428.9473
$ws.Range("I132").Value = 3209.2727
$ws.Range("J132").Value = 3731
$ws.Range("K132").Value = 9627.8181
$ws.Range("L132").Value = 11193
$ws.Range("M132").Value = -7097.8181
$ws.Range("N132").Value = -16253
$ws.Range("H136").Value = 1392.1428
$ws.Range("I136").Value = 1268.4615
$ws.Range("K136").Value = 3805.3845
$ws.Range("M136").Value = -1255.3845
